$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01514828764759746
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 9844.520545567508
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 9853.508716736955

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 5.553084769722144

$ws.Range("B4").Value = 0.127881588408715
$ws.Range("C4").Value = 0.04240448674262143
$ws.Range("D4").Value = 337.1190423067083
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 345.9495608678086

$ws.Range("B5").Value = 1.459612070389937
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.1575252929769615
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 11.945164432584

$ws.Range("B6").Value = 0.3048080303191223
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 3.900430680208489
$ws.Range("E6").Value = 8.660232485948974
$ws.Range("G6").Value = 14.53326577974471
